$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (week number + reporting date range)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 29   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# ---------------------------------------------------------------------------
# Helper: convert a cell that currently holds a number into a text cell
# (style 14 / General number format), matching the "0" / "***.*" markers
# used elsewhere on this sheet. We copy number-format only from a cell that
# is already styled that way (A14) so the underlying style index is reused
# instead of a brand-new one being allocated.
# ---------------------------------------------------------------------------
function Set-TextCell($addr, $text) {
    $ws.Range("A14").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = $false
    $ws.Range($addr).Value = $text
}

# Helper: convert a cell that currently holds text into a plain integer
# numeric cell (style 15 / "#,##0").
function Set-NumberCell($addr, $value) {
    $ws.Range($addr).Value = $value
    $ws.Range($addr).NumberFormat = "#,##0"
}

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
$ws.Range("N14").Value = -68.627450980392

# ---------------------------------------------------------------------------
# Row 15 - Rape (gains a week-to-date complaint, so the "0 -> ***.*" no-data
# markers move from column C/E to column D/E)
# ---------------------------------------------------------------------------
Set-NumberCell "C15" 1
Set-TextCell   "D15" "0"
Set-TextCell   "E15" "***.*"
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 28
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -30
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -70.212765957446

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 19
$ws.Range("G16").Value = 33
$ws.Range("H16").Value = -42.424242424242
$ws.Range("I16").Value = 307
$ws.Range("J16").Value = 240
$ws.Range("K16").Value = 27.916666666666
$ws.Range("L16").Value = 24.291497975708
$ws.Range("M16").Value = -21.882951653944
$ws.Range("N16").Value = -86.439929328621

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 54
$ws.Range("H17").Value = -40.740740740740
$ws.Range("I17").Value = 603
$ws.Range("J17").Value = 614
$ws.Range("K17").Value = -1.791530944625
$ws.Range("L17").Value = 14.638783269962
$ws.Range("M17").Value = 40.887850467289
$ws.Range("N17").Value = -48.680851063829

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 22.222222222222
$ws.Range("I18").Value = 212
$ws.Range("J18").Value = 151
$ws.Range("K18").Value = 40.397350993377
$ws.Range("L18").Value = 10.994764397905
$ws.Range("M18").Value = -47.654320987654
$ws.Range("N18").Value = -89.205702647657

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = -16.981132075471
$ws.Range("I19").Value = 643
$ws.Range("J19").Value = 454
$ws.Range("K19").Value = 41.629955947136
$ws.Range("L19").Value = 57.598039215686
$ws.Range("M19").Value = 24.371373307543
$ws.Range("N19").Value = -14.266666666666

# ---------------------------------------------------------------------------
# Row 20 - G.L.A. (week-to-date drops to zero, so column C becomes the
# "0" no-data marker)
# ---------------------------------------------------------------------------
Set-TextCell "C20" "0"
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -18.75
$ws.Range("J20").Value = 211
$ws.Range("K20").Value = 8.056872037914
$ws.Range("L20").Value = 38.181818181818
$ws.Range("M20").Value = -13.307984790874
$ws.Range("N20").Value = -88.253477588871

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = -11.111111111111
$ws.Range("F21").Value = 133
$ws.Range("G21").Value = 176
$ws.Range("H21").Value = -24.431818181818
$ws.Range("I21").Value = 2037
$ws.Range("J21").Value = 1711
$ws.Range("K21").Value = 19.053185271770
$ws.Range("L21").Value = 28.032683846637
$ws.Range("M21").Value = -1.068479844584
$ws.Range("N21").Value = -75.276125743415

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = 29.032258064516
$ws.Range("F24").Value = 120
$ws.Range("G24").Value = 88
$ws.Range("H24").Value = 36.363636363636
$ws.Range("I24").Value = 1066
$ws.Range("J24").Value = 871
$ws.Range("K24").Value = 22.388059701492
$ws.Range("L24").Value = 32.094175960347
$ws.Range("M24").Value = 18.708240534521

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -7.692307692307
$ws.Range("F25").Value = 56
$ws.Range("G25").Value = 58
$ws.Range("H25").Value = -3.448275862068
$ws.Range("I25").Value = 705
$ws.Range("J25").Value = 638
$ws.Range("K25").Value = 10.501567398119
$ws.Range("L25").Value = 29.357798165137
$ws.Range("M25").Value = -16.863207547169

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape* (gains a week-to-date complaint, no-data markers move
# from column C/E to column D/E)
# ---------------------------------------------------------------------------
Set-NumberCell "C26" 1
Set-TextCell   "D26" "0"
Set-TextCell   "E26" "***.*"
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 33.333333333333
$ws.Range("I26").Value = 47
$ws.Range("K26").Value = 20.512820512820
$ws.Range("L26").Value = -14.545454545454

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes (same pattern as row 26)
# ---------------------------------------------------------------------------
Set-NumberCell "C27" 1
Set-TextCell   "D27" "0"
Set-TextCell   "E27" "***.*"
$ws.Range("I27").Value = 73
$ws.Range("K27").Value = 23.728813559322
$ws.Range("L27").Value = 23.728813559322

# ---------------------------------------------------------------------------
# Row 28 - Shooting Vic. (week-to-date goes from zero to a reported value)
# ---------------------------------------------------------------------------
Set-NumberCell "C28" 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -60
$ws.Range("I28").Value = 48
$ws.Range("J28").Value = 58
$ws.Range("K28").Value = -17.241379310344
$ws.Range("L28").Value = -37.662337662337
$ws.Range("M28").Value = -39.240506329113
$ws.Range("N28").Value = -77.464788732394

# ---------------------------------------------------------------------------
# Row 29 - Shooting Inc. (same pattern as row 28)
# ---------------------------------------------------------------------------
Set-NumberCell "C29" 1
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = -60
$ws.Range("I29").Value = 41
$ws.Range("J29").Value = 52
$ws.Range("K29").Value = -21.153846153846
$ws.Range("L29").Value = -30.508474576271
$ws.Range("M29").Value = -39.705882352941
$ws.Range("N29").Value = -78.645833333333

# ---------------------------------------------------------------------------
# Row 30 - Hate Crimes
# ---------------------------------------------------------------------------
$ws.Range("L30").Value = -33.333333333333
